# B6-PowerPoint.pptx edit
#
# 1) Three tables (on the slides that hold them) get switched from the
#    custom "Table_0" style to a different built-in table style.
# 2) The presentation's theme colour scheme (currently the "Integral" /
#    "Red Violet" palette) is swapped for the stock "Office" palette.

$p = $ppt.ActivePresentation

# --- 1) Re-style the three tables -----------------------------------------
$newTableStyleId = "{9D1BB6A1-20E4-4E65-BB1C-846554A9BBDB}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# --- 2) Swap the theme colour scheme ---------------------------------------
# New ("Office") colours, in the standard dk1,lt1,dk2,lt2,accent1-6,hlink,
# folHlink order used by ThemeColorScheme.
$officeColors = @(
    0,            # dk1      000000
    16777215,     # lt1      FFFFFF
    6968388,      # dk2      44546A
    15132391,     # lt2      E7E6E6
    13998939,     # accent1  5B9BD5
    3243501,      # accent2  ED7D31
    10855845,     # accent3  A5A5A5
    49407,        # accent4  FFC000
    12874308,     # accent5  4472C4
    4697456,      # accent6  70AD47
    12673797,     # hlink    0563C1
    7491477       # folHlink 954F72
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($k = 1; $k -le $tcs.Count; $k++) {
    $tcs.Item($k).RGB = $officeColors[$k - 1]
}
